$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same formatting as the existing header cell (H1) to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(4, 5),
    @(6, 6),
    @(4, 5),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(3, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
